$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A30 needs to hold the literal text "2020-06-29" (not an auto-converted date serial).
# Assigning a date-shaped string straight into .Value/.Formula causes Excel's normal
# input parser to coerce it to a date (and stamp a date NumberFormat on the cell), so
# instead we compute it as a text formula and paste-special just the value back in,
# which keeps the cell as plain text with no style applied - exactly like the other
# date-label cells in column A.
$ws.Range("A30").Formula = '="2020-06-29"'
$ws.Range("A30").Copy()
$ws.Range("A30").PasteSpecial(-4163)

$ws.Range("B30").Value = 220657
$ws.Range("C30").Value = 279035
$ws.Range("D30").Value = 66910
$ws.Range("E30").Value = 27121
$ws.Range("F30").Value = 30.95
